$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11's rule-name label changes from "R40" to "1".
# Force the cell to keep storing a text value (not auto-convert the
# numeric-looking "1" into a number), matching the shared-string cell
# type ("t=s") used by the rest of the rule-name column.
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
